# Ajout de "DemanderMdp" + Modif page index
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Row 68 (CtrlDemanderMdp) - move "X" from "A faire" (B) to "En cours" (C),
# set Responsable (E) to "Killian" and Date de debut (F) to 27/09/2016
$ws.Range("B68").Value = $null
$ws.Range("C68").Value = "X"
$ws.Range("E68").Value = "Killian"
$ws.Range("F68").Value = (Get-Date -Year 2016 -Month 9 -Day 27 -Hour 0 -Minute 0 -Second 0).Date

# Row 69 (VueDemanderMdp) - same changes
$ws.Range("B69").Value = $null
$ws.Range("C69").Value = "X"
$ws.Range("E69").Value = "Killian"
$ws.Range("F69").Value = (Get-Date -Year 2016 -Month 9 -Day 27 -Hour 0 -Minute 0 -Second 0).Date

# Update view: select F69, then scroll so row 61 is the top visible row
$ws.Range("F69").Select()
$excel.ActiveWindow.ScrollRow = 61
